# Updates cryptos list values per the commit "Updated cryptos list on Fri Sep  6 23:42:18 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '53.851.21'
$ws.Cells.Item(2, 5).Value = '  -3.92%  '

$ws.Cells.Item(3, 4).Value = '2.221.02'
$ws.Cells.Item(3, 5).Value = '  -5.95%  '

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '0.996'
$ws.Cells.Item(4, 5).Value = '  -0.36%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '486.59'
$ws.Cells.Item(5, 5).Value = '  -3.12%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '124.91'
$ws.Cells.Item(6, 5).Value = '  -3.50%  '

$ws.Cells.Item(7, 5).Value = '  -0.28%  '

$ws.Cells.Item(8, 5).Value = '  -4.29%  '

$ws.Cells.Item(9, 4).Value = '2.217.06'
$ws.Cells.Item(9, 5).Value = '  -6.17%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.0923'
$ws.Cells.Item(10, 5).Value = '  -6.10%  '

$ws.Cells.Item(11, 5).Value = '  -1.32%  '

$ws.Cells.Item(13, 5).Value = '  -3.28%  '

$ws.Cells.Item(14, 4).Value = '2.587.26'
$ws.Cells.Item(14, 5).Value = '  -6.99%  '

$ws.Cells.Item(15, 5).Value = '  -1.45%  '

$ws.Cells.Item(16, 4).Value = '53.563.52'
$ws.Cells.Item(16, 5).Value = '  -4.29%  '

$ws.Cells.Item(17, 5).Value = '  -2.98%  '

$ws.Cells.Item(18, 4).Value = '2.215.95'
$ws.Cells.Item(18, 5).Value = '  -1.32%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '3.96'
$ws.Cells.Item(19, 5).Value = '  -1.44%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '9.56'
$ws.Cells.Item(20, 5).Value = '  -4.44%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '294.57'
$ws.Cells.Item(21, 5).Value = '  -4.25%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '6.14'
$ws.Cells.Item(22, 5).Value = '  -2.58%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '1.00'
$ws.Cells.Item(23, 5).Value = '  +0.12%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '62.76'
$ws.Cells.Item(24, 5).Value = '  -4.78%  '

$ws.Cells.Item(25, 5).Value = '  -0.27%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.364'
$ws.Cells.Item(26, 5).Value = '  -1.42%  '

$ws.Cells.Item(27, 4).Value = '2.299.83'
$ws.Cells.Item(27, 5).Value = '  -6.99%  '

$ws.Cells.Item(28, 5).Value = '  -1.34%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '7.01'
$ws.Cells.Item(29, 5).Value = '  -3.07%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '165.40'

$ws.Cells.Item(31, 5).Value = '  -3.80%  '

$ws.Cells.Item(33, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.999'
$ws.Cells.Item(33, 5).Value = '  +0.16%  '

$ws.Cells.Item(34, 2).Value = 'PEPE'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(34, 4).Value = '0.0₃0662'
$ws.Cells.Item(34, 5).Value = '  -6.69%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '5.69'
$ws.Cells.Item(35, 5).Value = '  -1.60%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '1.07'
$ws.Cells.Item(36, 5).Value = '  -1.11%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '17.29'
$ws.Cells.Item(37, 5).Value = '  -1.81%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '1.15'
$ws.Cells.Item(38, 5).Value = '  -1.50%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.828'
$ws.Cells.Item(39, 5).Value = '  +3.69%  '

$ws.Cells.Item(40, 2).Value = 'NEARProtocol'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '3.54'
$ws.Cells.Item(40, 5).Value = '  -4.72%  '

$ws.Cells.Item(41, 2).Value = 'OKB'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '35.81'
$ws.Cells.Item(41, 5).Value = '  -1.10%  '

$ws.Cells.Item(42, 5).Value = '  -1.02%  '

$ws.Cells.Item(43, 5).Value = '  -1.23%  '

$ws.Cells.Item(44, 2).Value = 'Aave'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '126.28'
$ws.Cells.Item(44, 5).Value = '  -2.21%  '

$ws.Cells.Item(45, 2).Value = 'Filecoin'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '3.27'
$ws.Cells.Item(45, 5).Value = '  -2.67%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '4.77'
$ws.Cells.Item(46, 5).Value = '  +2.11%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.0881'
$ws.Cells.Item(47, 5).Value = '  -2.39%  '

$ws.Cells.Item(48, 5).Value = '  -4.82%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '231.56'
$ws.Cells.Item(49, 5).Value = '  -2.65%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.0471'
$ws.Cells.Item(50, 5).Value = '  -2.09%  '

$ws.Cells.Item(51, 5).Value = '  -3.61%  '
